$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Asset Class labels in column B (and the stray D9 value keeps its
# single-space placeholder). Order below mirrors the order in which the
# distinct replacement strings were first introduced, so the rebuilt shared
# string table lines up with the target workbook.
$ws.Range("B3").Value = "U.S. High Yield Bonds"
$ws.Range("B2").Value = "Non-U.S. Bonds"
$ws.Range("B4").Value = "Non-U.S. Bonds"
$ws.Range("B5").Value = "Emerging Mkts Stks"
$ws.Range("B6").Value = "Real Estate"
$ws.Range("B7").Value = "U.S. Sm Cap Val"
$ws.Range("B8").Value = "U.S. Sm Cap Growth"
$ws.Range("B9").Value = "U.S. Sm Cap Val"
$ws.Range("B10").Value = "U.S. Sm Cap Growth"
$ws.Range("B11").Value = "U.S. Lg Cap Growth"
$ws.Range("B12").Value = "Commodities"
$ws.Range("B13").Value = "U.S. Lg Cap Val"
$ws.Range("B14").Value = "U.S. Investment Grade Bonds"
$ws.Range("B15").Value = "Emerging Mkts Stks"
$ws.Range("B16").Value = "U.S. Investment Grade Bonds"
$ws.Range("B17").Value = "U.S. Investment Grade Bonds"
$ws.Range("B18").Value = "Commodities"
$ws.Range("B19").Value = "U.S. Sm Cap Growth"
$ws.Range("B20").Value = "U.S. Sm Cap Val"
$ws.Range("B21").Value = "U.S. High Yield Bonds"
$ws.Range("B22").Value = "Commodities"
$ws.Range("B27").Value = "U.S. Mid Cap Growth"
$ws.Range("B23").Value = "Foreign Industrialzed Mkts Stocks"
$ws.Range("B24").Value = "U.S. Investment Grade Bonds"
$ws.Range("B25").Value = "U.S. Investment Grade Bonds"
$ws.Range("B26").Value = "U.S. Investment Grade Bonds"
$ws.Range("B28").Value = "U.S. Lg Cap Val"
$ws.Range("B29").Value = "U.S. Investment Grade Bonds"
$ws.Range("B30").Value = "U.S. Lg Cap Val"
$ws.Range("B31").Value = "U.S. Lg Cap Growth"
$ws.Range("B32").Value = "Emerging Mkts Stks"
$ws.Range("B33").Value = "Non-U.S. Bonds"
$ws.Range("B34").Value = "U.S. Sm Cap Val"
$ws.Range("B35").Value = "Foreign Industrialzed Mkts Stocks"

# Update the saved cursor/selection position shown in the sheet view.
$ws.Range("J23").Select() | Out-Null
